# Deploying to gh-pages from @ LinuxForHealth/alvearie-fhir-ig@80fa500adfae01c9a5dd7ef65e90accc96781b5c
# Re-brand the StructureDefinition workbook from IBM/Alvearie to LinuxForHealth:
#  - Metadata sheet: URL, Version, Date, Publisher
#  - Elements sheet: embedded extension URLs (Type(s) / Fixed Value columns)
#  - Elements sheet: drop the stray ele-1/ext-1 constraint text that had been
#    duplicated onto the root "Extension" row

$wb = $excel.ActiveWorkbook

# ---- Metadata sheet -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/sent-to-vendor"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# ---- Elements sheet ---------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# Extension.id / Element.id row no longer carries the ele-1/ext-1
# constraint text (that now only lives on the Extension.extension row).
$elements.Range("AI2").Value = ""

# communicationSentTime extension type
$elements.Range("J5").Value = "Extension {http://linuxforhealth.org/fhir/cdm/StructureDefinition/sent-time}" + [char]10

# communicationOutcome extension type
$elements.Range("J6").Value = "Extension {http://linuxforhealth.org/fhir/cdm/StructureDefinition/outcome}" + [char]10

# Extension.url fixed value
$elements.Range("Q7").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/sent-to-vendor"
